$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "код идентификатор PLU"
$ws.Range("B1").Value = "наименование анализируемых позиций"
$ws.Range("C1").Value = "данные по анализируемому критерию (продажи/оборот/прибыль)"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Товар 1"
$ws.Range("C2").Value = 100

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Товар 2"
$ws.Range("C3").Value = -50

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Товар 3"
$ws.Range("C4").Value = 35

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Товар 4"
$ws.Range("C5").Value = 20

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Товар 5"
$ws.Range("C6").Value = 5
